# Refresh odds figures (rows 2, 7-14) on Sheet1 to match the latest
# FlashScore snapshot pulled in this commit. Every assignment below
# corresponds to one changed <c>/<v> pair in the canonical OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 4.2
$ws.Range("J2").Value = 2.88
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("AD2").Value = 8
$ws.Range("AF2").Value = 17
$ws.Range("AO2").Value = 19

# Row 7
$ws.Range("H7").Value = 2.88
$ws.Range("I7").Value = 3.3
$ws.Range("S7").Value = 3.2
$ws.Range("T7").Value = 1.36
$ws.Range("W7").Value = 6.5
$ws.Range("X7").Value = 1.11
$ws.Range("AA7").Value = 2.5
$ws.Range("AB7").Value = 1.5
$ws.Range("AC7").Value = 5.5
$ws.Range("AE7").Value = 11
$ws.Range("AI7").Value = 5

# Row 8
$ws.Range("G8").Value = 2.7
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 2.55
$ws.Range("L8").Value = 3.4
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 9
$ws.Range("Q8").Value = 1.78
$ws.Range("R8").Value = 2.1
$ws.Range("S8").Value = 2.25
$ws.Range("T8").Value = 1.62
$ws.Range("W8").Value = 4
$ws.Range("X8").Value = 1.22
$ws.Range("AD8").Value = 13
$ws.Range("AF8").Value = 29
$ws.Range("AM8").Value = 351
$ws.Range("AP8").Value = 10

# Row 9
$ws.Range("G9").Value = 2.3
$ws.Range("H9").Value = 3.3
$ws.Range("J9").Value = 3.1
$ws.Range("K9").Value = 2
$ws.Range("Q9").Value = 1.78
$ws.Range("R9").Value = 2.1
$ws.Range("S9").Value = 2.35
$ws.Range("T9").Value = 1.57
$ws.Range("W9").Value = 4.33
$ws.Range("X9").Value = 1.2
$ws.Range("Y9").Value = 1.53
$ws.Range("Z9").Value = 2.38
$ws.Range("AA9").Value = 2
$ws.Range("AB9").Value = 1.75
$ws.Range("AF9").Value = 21
$ws.Range("AG9").Value = 21
$ws.Range("AH9").Value = 34
$ws.Range("AI9").Value = 8
$ws.Range("AM9").Value = 451
$ws.Range("AP9").Value = 11

# Row 10
$ws.Range("G10").Value = 3.5
$ws.Range("I10").Value = 2.15
$ws.Range("J10").Value = 4.5
$ws.Range("Q10").Value = 1.93
$ws.Range("R10").Value = 1.93
$ws.Range("AC10").Value = 8
$ws.Range("AD10").Value = 17
$ws.Range("AJ10").Value = 6.5
$ws.Range("AL10").Value = 81
$ws.Range("AO10").Value = 9
$ws.Range("AP10").Value = 9.5
$ws.Range("AQ10").Value = 19

# Row 11
$ws.Range("Q11").Value = 1.8
$ws.Range("R11").Value = 2.05
$ws.Range("S11").Value = 2.35
$ws.Range("T11").Value = 1.57

# Row 12
$ws.Range("G12").Value = 2.15
$ws.Range("I12").Value = 3.3
$ws.Range("J12").Value = 2.88
$ws.Range("N12").Value = 9
$ws.Range("AD12").Value = 9.5
$ws.Range("AF12").Value = 19
$ws.Range("AN12").Value = 9
$ws.Range("AQ12").Value = 41

# Row 13
$ws.Range("G13").Value = 1.38
$ws.Range("I13").Value = 8
$ws.Range("L13").Value = 7.5
$ws.Range("S13").Value = 1.75
$ws.Range("T13").Value = 2.05
$ws.Range("W13").Value = 2.75
$ws.Range("X13").Value = 1.4
$ws.Range("AC13").Value = 6.5
$ws.Range("AM13").Value = 401
$ws.Range("AQ13").Value = 101

# Row 14
$ws.Range("L14").Value = 3.75
$ws.Range("O14").Value = 1.4
$ws.Range("P14").Value = 2.75
$ws.Range("S14").Value = 2.3
$ws.Range("T14").Value = 1.6
$ws.Range("W14").Value = 4.33
$ws.Range("X14").Value = 1.2
$ws.Range("Y14").Value = 1.5
$ws.Range("Z14").Value = 2.5
$ws.Range("AA14").Value = 1.91
$ws.Range("AB14").Value = 1.8
$ws.Range("AC14").Value = 7
$ws.Range("AD14").Value = 10
$ws.Range("AI14").Value = 7.5
$ws.Range("AK14").Value = 15
$ws.Range("AM14").Value = 900
